$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 24.50000000000039
$ws.Range("H2").Value = [double]"1.215018358002907e-16"
$ws.Range("I2").Value = [double]"8.199634704553205e-08"
$ws.Range("K2").Value = 51.57029525989654
$ws.Range("L2").Value = "[44.49770505942, 58.64288546037308]"
$ws.Range("O2").Value = 1.742184514603349
$ws.Range("P2").Value = "[1.591237119836273, 1.8931319093704255]"
$ws.Range("S2").Value = 58.2733107637474
$ws.Range("T2").Value = "[53.88321426151235, 62.66340726598246]"
$ws.Range("W2").Value = 17.70670670670699
$ws.Range("X2").Value = 17.11811811811839
$ws.Range("Y2").Value = 18.29529529529558

# Row 3
$ws.Range("E3").Value = 23.75000000000027
$ws.Range("H3").Value = [double]"1.215018358002907e-16"
$ws.Range("I3").Value = [double]"6.493473636570712e-08"
$ws.Range("K3").Value = 51.18291270478535
$ws.Range("L3").Value = "[42.27254757518273, 60.093277834387976]"
$ws.Range("O3").Value = 0.3962369112635775
$ws.Range("P3").Value = "[0.22013161736865428, 0.5723422051585008]"
$ws.Range("Q3").Value = [double]"1.424525017035272e-05"
$ws.Range("R3").Value = [double]"1.424525017035272e-05"
$ws.Range("S3").Value = 54.97536347462677
$ws.Range("T3").Value = "[50.150503532446415, 59.800223416807135]"
$ws.Range("W3").Value = 22.25225225225251
$ws.Range("X3").Value = 21.58658658658683
$ws.Range("Y3").Value = 22.91791791791818
